# Control brazo - add Hoja2 (pwm -> angle in radians table) and
# underline the D3 header-ish cell on Hoja1, matching commit
# "Control brazo actualizado y clase de variables de estados"

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Hoja1 tweaks -----------------------------------------------------
# Underline D3 (new style: font with <u/>) and move the active selection
# to D3 (was M16).
$ws1.Range("D3").Font.Underline = $true
$ws1.Range("D3").Select() | Out-Null

# --- New sheet "Hoja2" -------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Hoja2"

$ws2.Range("A1").Value = "Pwm"
$ws2.Range("B1").Value = "Ángulo (Rad)"
$ws2.Range("A1:B1").Font.Underline = $true

$pwmValues = @(0,50,100,150,200,250,300,350,400,450,500,550,600,650,700,750,800,850,900,950,1000)
for ($i = 0; $i -lt $pwmValues.Length; $i++) {
    $row = $i + 2
    $ws2.Cells.Item($row, 1).Value = $pwmValues[$i]
    $ws2.Cells.Item($row, 2).Formula = "=A$row*PI()/180"
}

$ws2.Range("B6").Select() | Out-Null

$ws1.Select() | Out-Null
